$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.214.99"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.687.75"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +12.44%  "
$ws.Range("E9").Value = "  +4.64%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.926.06"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.684.10"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "27.221.02"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "0.0₃0743"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "1.550.98"
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").Value = "1.834.44"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0110"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("E50").Value = "  +9.09%  "
$ws.Range("E51").Value = "  +1.65%  "

# Cells whose new value looks numeric; force them to remain text
# to match the source data (which stores all prices as text strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.788"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.45"
$ws.Range("D50").Style = "Normal"

Write-Output "Updated cryptos list"